# Generate Report for Archive
#
# The localization "Status" for the tracked file moved from
# "Ready for handoff" to "In Translation". That status value is shown
# on the Overview sheet (once per target locale column: zh-cn, de-de)
# and on each per-locale detail sheet (zh-cn, de-de) in their own
# "Status" column. Updating the text makes the status column narrower,
# so its column width shrinks to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de), row 2 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

# --- zh-cn detail sheet: column C ("Status"), row 2 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- de-de detail sheet: column C ("Status"), row 2 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C1").ColumnWidth = 12.5
